$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: fill column A (subject names) for rows 3-13 first.
$ws.Range("A3").Value = "PPKN"
$ws.Range("A4").Value = "Bahasa Indonesia"
$ws.Range("A5").Value = "Matematika"
$ws.Range("A6").Value = "IPA"
$ws.Range("A7").Value = "IPS"
$ws.Range("A8").Value = "Bahasa Inggris"
$ws.Range("A9").Value = "Seni Budaya"
$ws.Range("A10").Value = "Penjas Orkes"
$ws.Range("A11").Value = "Prakarya"
$ws.Range("A12").Value = "Bahasa Jawa (muatan lokal)"
$ws.Range("A13").Value = "Pengembangan Diri(BK)"

# Step 2: fill row 2 (Religious education row, inserted above the rest).
$ws.Range("A2").Value = "Pend. Agama & Budi Pekerti"
$ws.Range("B2").Value = "PA"

# Step 3: fill column B (subject codes) for rows 3-13.
$ws.Range("B3").Value = "PPKN"
$ws.Range("B4").Value = "BIN"
$ws.Range("B5").Value = "MTK"
$ws.Range("B6").Value = "IPA"
$ws.Range("B7").Value = "IPS"
$ws.Range("B8").Value = "BIG"
$ws.Range("B9").Value = "SB"
$ws.Range("B10").Value = "OR"
$ws.Range("B11").Value = "PK"
$ws.Range("B12").Value = "BJ"
$ws.Range("B13").Value = "BK"

# Update the view selection to match the target state.
$ws.Range("B12").Select()

# Apply page setup (paper size 70 = "Japanese Postcard" in Excel's enumeration,
# matches the target pageSetup element written during printing/print-preview).
$ws.PageSetup.PaperSize = 70
$ws.PageSetup.Orientation = 1
